# updated legacy GSC export data
# The "Chart" sheet (sheet1) holds a rolling window of dates (col A),
# a constant 0 column (col B, "Non-HTTPS URLs") and the HTTPS URL counts
# (col C). The export window rolled forward by two days: the two oldest
# dates (2025-09-10, 2025-09-11) were dropped and four new trailing
# dates (2025-12-07 .. 2025-12-10) were appended with a count of 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the two oldest data rows (2025-09-10, 2025-09-11) - this shifts
# every remaining row up by two, carrying its existing date text/value
# along with it (no retyping, so no accidental date-number conversion).
$ws.Rows("2:3").Delete()

# After the delete, the last existing data row is row 87 (2025-12-06).
# Append four new rows continuing the daily sequence, each with
# Non-HTTPS URLs = 0 and HTTPS URLs = 0.
$newDates = @("2025-12-07", "2025-12-08", "2025-12-09", "2025-12-10")
$firstNewRow = 88

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $firstNewRow + $i
    # Leading apostrophe forces literal text so the ISO date string is
    # stored as-is instead of being parsed into a date serial number.
    $ws.Cells.Item($r, 1).Value = "'" + $newDates[$i]
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
}

# Writing a quoted literal marks the cell with a quote-prefix style;
# re-apply the plain formatting of the preceding (unmodified) row so the
# new cells match the rest of the column (General, no special style).
$lastOldRow = $firstNewRow - 1
$ws.Range("A" + $lastOldRow + ":C" + $lastOldRow).Copy()
$ws.Range("A" + $firstNewRow + ":C" + ($firstNewRow + $newDates.Length - 1)).PasteSpecial(-4122)
$excel.CutCopyMode = $false
